$d = $word.ActiveDocument

# Title (appears twice: H1 heading and bold run near the end) - replace all occurrences
$d.Content.Find.Execute("Play Atlantean GigaRise for Free - Game Review", $true, $false, $false, $false, $false, $true, 1, $false, "Play Atlantean GigaRise", 2)

# "What we like" bullet list
$d.Content.Find.Execute("Visually stunning graphics and animations.", $true, $false, $false, $false, $false, $true, 1, $false, "Gameplay mechanics and features", 2)
$d.Content.Find.Execute("Up to 294,912 ways to win.", $true, $false, $false, $false, $false, $true, 1, $false, "High-level graphics and visually beautiful underwater scene", 2)
$d.Content.Find.Execute("Automatic spins up to 1,000.", $true, $false, $false, $false, $false, $true, 1, $false, "Music accompanies players throughout the gaming session", 2)
$d.Content.Find.Execute("Fair RTP value of 96%.", $true, $false, $false, $false, $false, $true, 1, $false, "Special symbols with free spins for added excitement", 2)

# "What we don't like" bullet list
$d.Content.Find.Execute("High volatility may not suit new players.", $true, $false, $false, $false, $false, $true, 1, $false, "High volatility may not be suitable for inexperienced players", 2)
$d.Content.Find.Execute("Free spins feature is challenging to trigger.", $true, $false, $false, $false, $false, $true, 1, $false, "Limited betting limits with a maximum of €40", 2)

# Meta description paragraph
$d.Content.Find.Execute("Read our review of Atlantean GigaRise, a visually stunning slot game with up to 294,912 ways to win. Play it for free and learn about its features and betting limits.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Atlantean GigaRise to play this exciting slot game for free.", 2)
